$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "all-tabs" StatQuery text replacing the old one that counted files/
# samples/cases/studies. The new query additionally counts programs and
# distinguishes case-files from study-files.
$newStatQuery = "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (diag:diagnosis)-->(c)`nOPTIONAL MATCH (f:file)-[*]->(c)`nOPTIONAL MATCH (sf:file)-->(s)`nWITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p`nWHERE demo.breed IN ['Belgian Malinois']`nRETURN  `n    count(distinct p) AS Programs,`n    count(distinct s) AS Studies,`n    count(distinct c) AS Cases,`n    count(distinct samp) AS Samples,`n    count(distinct f) AS ``Case Files``,`n    count(distinct sf) AS ``Study Files``"

# Replace the StatQuery column (C) on the Cases/Samples/Files rows with the
# new query text.
$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Update the saved cursor position/selection to B5 (and drop the scrolled
# topLeftCell, i.e. scroll the view back to the top).
$ws.Range("B5").Select()
